$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (workbook.xml sheet name change)
$ws.Name = "BOM_Board1_PCB1_2023-09-22"

# Update BOM row 5 (part No. 4) values - adapter PCB connector changed
$ws.Range("C5").Value = "HDGC2001WR-4P"
$ws.Range("G5").Value = "HDGC2001WR-4P"
$ws.Range("D5").Value = "VE.DIRECT1"
$ws.Range("E5").Value = "CONN-TH_4P-P2.00_HDGC_HDGC2001WR-4P"
$ws.Range("H5").Value = "HDGC(华德共创)"
$ws.Range("I5").Value = "C5175241"
